$wb = $excel.ActiveWorkbook

# Update metadata values on the "Metadata" sheet
$meta = $wb.Worksheets("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# Update the "descendent-of" value on the "Include from FSIII" sheet
$inc1 = $wb.Worksheets("Include from FSIII")
$inc1.Range("C2").Value = "G1"

# Remove the two trailing "Include from FSIII" sheets (revert to 1.1.0 shape)
$wb.Worksheets("Include from FSIII 4").Delete()
$wb.Worksheets("Include from FSIII 3").Delete()

# Restore original active-sheet selection (Metadata, first tab)
$meta.Activate()
